$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the floating textbox's current vertical position so we can keep
# it visually anchored to the same place after the row shift below (the
# shape otherwise stays pinned to its absolute position and does not
# automatically follow the deleted row).
$shape = $ws.Shapes.Item(1)
$shapeTop = $shape.Top
$row1Height = $ws.Rows.Item(1).RowHeight

# Delete the first (blank) row, shifting the header and data rows up by one.
$ws.Rows.Item(1).Delete()

# Re-anchor the textbox so it moves up along with the rows above it.
$shape2 = $ws.Shapes.Item(1)
$shape2.Top = $shapeTop - $row1Height

# Update selection to match target state (A5 selected).
$ws.Range("A5").Select()
